# Update "want to go" counts (column F) for a handful of events.
# Same rows/values need to change on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    8  = 499
    10 = 1927
    11 = 58
    13 = 4006
    19 = 50
    20 = 2728
    33 = 1576
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
